$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 3 with new credentials
$ws.Range("A3").Value = "mhsZu871"
$ws.Range("B3").Value = "pZgbiz"

# Remove rows 4 through 6 (old extra credential rows)
$ws.Range("A4:B6").EntireRow.Delete()

# Update the active selection to A6
$ws.Range("A6").Select()
